# Auto-generated edit script: updates cryptos list values per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force these D-column cells to Text format first so values such as "1.00"
# or "58.90" are preserved verbatim instead of being normalised to numbers.
$textCells = @(
    "D2", "D3", "D4", "D5", "D6", "D8", "D9", "D10", "D11", "D12",
    "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23",
    "D24", "D25", "D27", "D28", "D29", "D30", "D32", "D33", "D34", "D35",
    "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46",
    "D47", "D48", "D49", "D50", "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values row by row.
# Row 2
$ws.Range("D2").Value = "55.363.37"

# Row 3
$ws.Range("D3").Value = "2.915.98"
$ws.Range("E3").Value = "  -9.87%  "

# Row 4
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").Value = "468.63"
$ws.Range("E5").Value = "  -13.01%  "

# Row 6
$ws.Range("D6").Value = "123.85"
$ws.Range("E6").Value = "  -9.42%  "

# Row 7
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("D8").Value = "2.915.32"
$ws.Range("E8").Value = "  -9.89%  "

# Row 9
$ws.Range("D9").Value = "0.401"
$ws.Range("E9").Value = "  -12.70%  "

# Row 10
$ws.Range("D10").Value = "6.54"
$ws.Range("E10").Value = "  -14.17%  "

# Row 11
$ws.Range("D11").Value = "0.0946"
$ws.Range("E11").Value = "  -17.92%  "

# Row 12
$ws.Range("D12").Value = "0.327"
$ws.Range("E12").Value = "  -17.52%  "

# Row 13
$ws.Range("E13").Value = "  -3.66%  "

# Row 14
$ws.Range("D14").Value = "3.418.10"
$ws.Range("E14").Value = "  -9.78%  "

# Row 15
$ws.Range("D15").Value = "22.37"
$ws.Range("E15").Value = "  -14.61%  "

# Row 16
$ws.Range("D16").Value = "55.338.24"
$ws.Range("E16").Value = "  -6.31%  "

# Row 17
$ws.Range("D17").Value = "2.922.44"
$ws.Range("E17").Value = "  -9.63%  "

# Row 18
$ws.Range("D18").Value = "0.0000131"
$ws.Range("E18").Value = "  -17.48%  "

# Row 19
$ws.Range("D19").Value = "5.06"
$ws.Range("E19").Value = "  -14.00%  "

# Row 20
$ws.Range("D20").Value = "11.41"
$ws.Range("E20").Value = "  -13.66%  "

# Row 21
$ws.Range("D21").Value = "6.97"
$ws.Range("E21").Value = "  -15.86%  "

# Row 22
$ws.Range("D22").Value = "305.34"
$ws.Range("E22").Value = "  -15.56%  "

# Row 23
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.08%  "

# Row 24
$ws.Range("D24").Value = "0.442"
$ws.Range("E24").Value = "  -15.02%  "

# Row 25
$ws.Range("D25").Value = "58.90"
$ws.Range("E25").Value = "  -16.72%  "

# Row 26
$ws.Range("E26").Value = "  +0.49%  "

# Row 27
$ws.Range("D27").Value = "0.155"
$ws.Range("E27").Value = "  -9.11%  "

# Row 28
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  -0.13%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0789"
$ws.Range("E29").Value = "  -19.33%  "

# Row 30
$ws.Range("D30").Value = "5.83"
$ws.Range("E30").Value = "  -17.25%  "

# Row 31
$ws.Range("E31").Value = "  -9.77%  "

# Row 32
$ws.Range("D32").Value = "18.87"
$ws.Range("E32").Value = "  -14.40%  "

# Row 33
$ws.Range("D33").Value = "5.97"
$ws.Range("E33").Value = "  -15.20%  "

# Row 34
$ws.Range("D34").Value = "1.58"
$ws.Range("E34").Value = "  -18.62%  "

# Row 35
$ws.Range("D35").Value = "144.84"
$ws.Range("E35").Value = "  -11.63%  "

# Row 36
$ws.Range("E36").Value = "  -15.89%  "

# Row 37
$ws.Range("D37").Value = "5.34"
$ws.Range("E37").Value = "  -16.46%  "

# Row 38
$ws.Range("D38").Value = "1.20"
$ws.Range("E38").Value = "  -16.09%  "

# Row 39
$ws.Range("D39").Value = "2.944.91"
$ws.Range("E39").Value = "  -9.84%  "

# Row 40
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.02%  "

# Row 41
$ws.Range("D41").Value = "0.0603"
$ws.Range("E41").Value = "  -14.86%  "

# Row 42
$ws.Range("D42").Value = "21.19"
$ws.Range("E42").Value = "  -19.05%  "

# Row 43
$ws.Range("D43").Value = "34.91"
$ws.Range("E43").Value = "  -15.18%  "

# Row 44
$ws.Range("D44").Value = "0.954"
$ws.Range("E44").Value = "  -12.91%  "

# Row 45
$ws.Range("D45").Value = "0.597"
$ws.Range("E45").Value = "  -16.88%  "

# Row 46
$ws.Range("D46").Value = "3.36"
$ws.Range("E46").Value = "  -16.24%  "

# Row 47
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").Value = "1.29"
$ws.Range("E47").Value = "  -14.75%  "

# Row 48
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "2.027.05"
$ws.Range("E48").Value = "  -11.69%  "

# Row 49
$ws.Range("D49").Value = "5.23"
$ws.Range("E49").Value = "  -16.73%  "

# Row 50
$ws.Range("D50").Value = "17.38"
$ws.Range("E50").Value = "  -16.13%  "

# Row 51
$ws.Range("D51").Value = "0.0207"
$ws.Range("E51").Value = "  -14.40%  "

